$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2025-05-30 Friday" "2025-05-31 Saturday"

Replace-Text "644×5=" "614×6="
Replace-Text "915×8=" "713×3="
Replace-Text "647×5=" "133×8="
Replace-Text "307×8=" "342×2="
Replace-Text "503×8=" "781×9="
Replace-Text "883×8=" "124×7="
Replace-Text "406×9=" "723×3="
Replace-Text "417×5=" "309×6="
Replace-Text "465×7=" "249×3="
Replace-Text "907×9=" "898×3="
Replace-Text "229×3=" "749×5="
Replace-Text "822×5=" "271×8="
Replace-Text "304×5=" "664×6="
Replace-Text "299×3=" "884×3="
Replace-Text "899×5=" "404×2="
Replace-Text "547×4=" "540×7="
Replace-Text "803×3=" "801×3="
Replace-Text "557×6=" "796×6="
Replace-Text "569×8=" "935×2="
Replace-Text "534×8=" "498×2="
Replace-Text "380×4=" "301×4="
Replace-Text "488×7=" "164×9="
Replace-Text "390×4=" "838×8="
Replace-Text "453×8=" "671×9="
Replace-Text "252×9=" "350×7="
